$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OCT-2020")
$ws.Activate()

# Fill in Application (C) and Comments (D) columns for rows 18-20 and 23-25
$ws.Range("C18").Value = "QMVAR and Haayai"
$ws.Range("D18").Value = "Qmvar issues going on, for hayaai for token system new table created"

$ws.Range("C19").Value = "QMVAR and Haayai"
$ws.Range("D19").Value = "Hayaai app datas issues solved and Qmvar issues going on"

$ws.Range("C20").Value = "QMVAR "
$ws.Range("D20").Value = "In Qmvar upload summary details close button issue going on."

$ws.Range("C23").Value = " Mujistore and Hayaai"
$ws.Range("D23").Value = "Fixing 3 issues from mujistore and support database issues for Hayaai app"

$ws.Range("C24").Value = " Mujistore and Hayaai"
$ws.Range("D24").Value = "Fixing 2 issues from mujistore and support database issues for Hayaai app"

$ws.Range("C25").Value = " Mujistore and Hayaai"
$ws.Range("D25").Value = "Deployment given for Mujistore and support database work for Hayaai app"

# Update the saved view state (scroll position + active selection)
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("C25").Select()
